# The source data's "HS" (Highest Score) column (G) stored values such as
# "39*" as text (the trailing "*" denotes "not out"). The authoritative
# edit strips the "*" and re-saves the column as a plain number for every
# row where it was previously a starred text value. This also causes the
# now-unused "NN*" shared strings to be dropped from the shared string
# table on save (and all later string indices to shift down), which is
# purely a side effect of the engine's string-table compaction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Val = 39 }
    @{ Row = 6;  Val = 0 }
    @{ Row = 7;  Val = 58 }
    @{ Row = 15; Val = 76 }
    @{ Row = 17; Val = 100 }
    @{ Row = 20; Val = 59 }
    @{ Row = 22; Val = 43 }
    @{ Row = 23; Val = 35 }
    @{ Row = 25; Val = 38 }
    @{ Row = 27; Val = 15 }
    @{ Row = 29; Val = 57 }
    @{ Row = 31; Val = 124 }
    @{ Row = 32; Val = 63 }
    @{ Row = 33; Val = 75 }
    @{ Row = 34; Val = 64 }
    @{ Row = 35; Val = 66 }
    @{ Row = 36; Val = 68 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Val
}

# Match the author's final selection / scroll position in the saved view.
$ws.Range("G2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1

$wb.Save()
